$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "Changes"
$ws.Range("N1").Value = "Notes"

$ws.Range("A3").Value = 201911271207
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = "MG laptop"
$ws.Range("D3").Value = "TDB2 persistent"
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 500
$ws.Range("G3").Value = 0.7934
$ws.Range("H3").Value = 2.0379
$ws.Range("I3").Value = 622.394
$ws.Range("J3").Value = 2892.06
$ws.Range("K3").Value = 1.24479
$ws.Range("L3").Value = 1.22952
$ws.Range("M3").Value = "Graph name cache enabled"
$ws.Range("N3").Value = "~ Halve runtime and double throughput"

$ws.Range("L4").Select()
